$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1397
$ws.Range("F4").Value = 13427
$ws.Range("F5").Value = 778
$ws.Range("F13").Value = 21880
$ws.Range("F14").Value = 545
$ws.Range("F15").Value = 224
$ws.Range("F16").Value = 515
$ws.Range("F17").Value = 137
$ws.Range("F18").Value = 376
$ws.Range("F21").Value = 165
$ws.Range("F22").Value = 144
$ws.Range("F24").Value = 236
$ws.Range("F27").Value = 1364
$ws.Range("F28").Value = 67
$ws.Range("F29").Value = 381
$ws.Range("F31").Value = 108

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 304
$ws.Range("F5").Value = 195
$ws.Range("F8").Value = 17
$ws.Range("F10").Value = 89
$ws.Range("F11").Value = 89
$ws.Range("F12").Value = 391
$ws.Range("D18").Value = "广州大道中1229号 广东艺术剧院"
$ws.Range("F18").Value = 20

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 4478
$ws.Range("F4").Value = 108

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 1397
$ws.Range("F6").Value = 13427
$ws.Range("F7").Value = 304
$ws.Range("F8").Value = 778
$ws.Range("F9").Value = 4478
$ws.Range("F15").Value = 108
$ws.Range("F16").Value = 21889
$ws.Range("F17").Value = 545
$ws.Range("F19").Value = 224
$ws.Range("F20").Value = 195
$ws.Range("F21").Value = 195
$ws.Range("F22").Value = 515
$ws.Range("F25").Value = 137
$ws.Range("F26").Value = 17
$ws.Range("F28").Value = 89
$ws.Range("F29").Value = 391
$ws.Range("F30").Value = 376
$ws.Range("F33").Value = 165
$ws.Range("F34").Value = 144
$ws.Range("F37").Value = 236
$ws.Range("F42").Value = 1364
$ws.Range("F43").Value = 67
$ws.Range("F45").Value = 381
$ws.Range("F47").Value = 108
$ws.Range("D48").Value = "广州大道中1229号 广东艺术剧院"
$ws.Range("F48").Value = 20
